# Generate Report for Handoff
# Adds two new localization entries (6ab99848-...md and 6e44e4f3-...md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Overview sheet - rows 6 & 7
# ---------------------------------------------------------------------------
$wsOverview.Range("A6").Value = "6ab99848-e29c-4db5-af0b-f331d60c43cf.md"
$wsOverview.Range("B6").Value = "e2e\6ab99848-e29c-4db5-af0b-f331d60c43cf.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-25 08:44:34"

$wsOverview.Range("A7").Value = "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md"
$wsOverview.Range("B7").Value = "e2e\6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-25 08:44:34"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab99848-e29c-4db5-af0b-f331d60c43cf/e2e/6ab99848-e29c-4db5-af0b-f331d60c43cf.md", "", "", "e2e\6ab99848-e29c-4db5-af0b-f331d60c43cf.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6e44e4f3-8fea-4c11-b086-0346aaaa2b23/e2e/6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md", "", "", "e2e\6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md")

# Grow the "Overview" table / autofilter / dimension to A1:G7
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# 2. zh-cn sheet - rows 6 & 7
# ---------------------------------------------------------------------------
$wsZhCn.Range("A6").Value = "6ab99848-e29c-4db5-af0b-f331d60c43cf.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "6ab99848-e29c-4db5-af0b-f331d60c43cf.45557f6b6bc5fa0279f64084fd49eb4a56c5ee4e.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-25 08:44:29"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("O6").Value = "False"

$wsZhCn.Range("A7").Value = "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.1299da7a356c32e3231dc87c7fe1891ff6bca477.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-25 08:44:29"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6ab99848-e29c-4db5-af0b-f331d60c43cf/e2e/6ab99848-e29c-4db5-af0b-f331d60c43cf.md", "", "", "6ab99848-e29c-4db5-af0b-f331d60c43cf.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6e44e4f3-8fea-4c11-b086-0346aaaa2b23/e2e/6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md", "", "", "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------------
# 3. de-de sheet - rows 6 & 7
# ---------------------------------------------------------------------------
$wsDeDe.Range("A6").Value = "6ab99848-e29c-4db5-af0b-f331d60c43cf.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "6ab99848-e29c-4db5-af0b-f331d60c43cf.45557f6b6bc5fa0279f64084fd49eb4a56c5ee4e.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-25 08:44:34"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("O6").Value = "False"

$wsDeDe.Range("A7").Value = "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.1299da7a356c32e3231dc87c7fe1891ff6bca477.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-25 08:44:34"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6ab99848-e29c-4db5-af0b-f331d60c43cf/e2e/6ab99848-e29c-4db5-af0b-f331d60c43cf.md", "", "", "6ab99848-e29c-4db5-af0b-f331d60c43cf.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6e44e4f3-8fea-4c11-b086-0346aaaa2b23/e2e/6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md", "", "", "6e44e4f3-8fea-4c11-b086-0346aaaa2b23.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P7"))

Write-Host "Report generated for handback."
